$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: find the row number of the data row whose "Conta" (column A)
# equals the given account number.
function Get-AccountRow($account) {
    $cell = $ws.Columns(1).Find($account)
    return $cell.Row
}

# 1. Remove the SURAMA row (004205505 / SURAMA / 21345.59).
$suramaRow = Get-AccountRow "004205505"
$ws.Rows($suramaRow).Delete()

# 2. Insert a new row for ALBERTO (004480970 / ALBERTO / 7633.09) directly
#    above MONICA (005170415).
$monicaRow = Get-AccountRow "005170415"
$ws.Rows($monicaRow).Insert()
$albertoRow = $monicaRow
$ws.Cells.Item($albertoRow, 1).NumberFormat = "@"
$ws.Cells.Item($albertoRow, 1).Value = "004480970"
$ws.Cells.Item($albertoRow, 2).Value = "ALBERTO"
$ws.Cells.Item($albertoRow, 3).Value = 7633.09

# 3. Remove the block of negative-balance rows that were dropped from the
#    bottom of the sheet (everyone below THIAGO's -41892.08 balance,
#    except THIAGO himself, plus a run of larger negative balances just
#    above him).
$accountsToRemove = @(
    "004630773", # NABOR        -6444.39
    "004948033", # GUILHERME    -7518.44
    "005002457", # ROSANGELA    -7518.44
    "005255637", # PATRICIA     -8234.54
    "004752534", # CARLOS       -8473.29
    "005105970", # VERA         -8711.91
    "004751154", # CATARINE     -8831.33
    "004453132", # BRUNO       -16231.14
    "004499920", # FABIANO     -21243.78
    "004556974", # KELLY       -21721.31
    "004482102", # NATALIA     -22079.35
    "005009992", # ALINE       -29240.45
    "004940560", # CRISTIANO   -40076.32
    "004582648", # ANUAR       -51798.26
    "004450760"  # SILVIO      -59317.56
)

# Resolve every row number up-front, then delete from the bottom row
# upward so earlier deletions never invalidate a still-pending row number.
$rowsToDelete = $accountsToRemove | ForEach-Object { Get-AccountRow $_ }
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
